# Add a "betreuer" (supervisor) column to the client_data sheet and fill it
# in with supervisor names, associating clients 1+2 with "Betreuer 1" and
# clients 3+4 with "Betreuer 2". Also clears the stray "zeichnungssumme"
# value that had been left for client 3's row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("client_data")

# Insert a new, empty column before column B; everything from B onward
# (headers, data, column widths) shifts one column to the right.
$ws.Columns("B").Insert()

# Header for the newly inserted column.
$ws.Range("B1").Value = "betreuer"

# Supervisor assignment per client row.
$ws.Range("B2").Value = "Betreuer 1"
$ws.Range("B3").Value = "Betreuer 1"
$ws.Range("B4").Value = "Betreuer 2"
$ws.Range("B5").Value = "Betreuer 2"

# Row 4 (client id 3) had an errant "zeichnungssumme" (subscription amount)
# value; that value (now shifted into column N) is removed.
$ws.Range("N4").ClearContents() | Out-Null

# Move the active selection, matching the saved workbook state.
$ws.Range("B6").Select() | Out-Null
